$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.031.54"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.269.13"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.656"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "233.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +6.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0984"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "2.606.57"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").Value = "2.266.75"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "43.946.86"
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +18.42%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.128"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("E34").Value = "  +7.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0686"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -4.72%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0954"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("D50").Value = "1.454.79"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.25%  "
